$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text (cell A1) with the new conversion rates ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $hoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 1.56 = 5715.82 pesos"), "✅ 1000 Bs = 1.64 = 6033.22 pesos"
$text = $text -replace [regex]::Escape("✅ 5715.82 pesos = 1.55 = 898.42 Bs"), "✅ 6033.22 pesos = 1.64 = 909.56 Bs"
$cell.Value2 = $text

# --- Update the "tasas" sheet numeric cells ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 608
$tasas.Range("O10").Value = 3668.2
$tasas.Range("N12").Value = 3688
$tasas.Range("O12").Value = 556
